$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Latitude (B) and Longitude (C) values for rows 2-10
$ws.Range("B2").Value = 51.0339991
$ws.Range("C2").Value = -2.94901191602119

$ws.Range("B3").Value = 51.50002095
$ws.Range("C3").Value = -0.192442035662122

$ws.Range("B4").Value = 51.4510190999999
$ws.Range("C4").Value = -0.993491440497515

$ws.Range("B5").Value = 52.219977
$ws.Range("C5").Value = 0.487578

$ws.Range("B6").Value = 52.45124
$ws.Range("C6").Value = -1.937937

$ws.Range("B7").Value = 54.0680924499999
$ws.Range("C7").Value = -2.68525312594421

$ws.Range("B8").Value = 56.7861112
$ws.Range("C8").Value = -4.1140518

$ws.Range("B9").Value = 52.2928116
$ws.Range("C9").Value = -3.73893

$ws.Range("B10").Value = 54.5859836
$ws.Range("C10").Value = -6.9591554

# Update the active selection to C1
$ws.Range("C1").Select()
